# edit.ps1 -- reproduce the OOXML diff against before.pptx
#
# Summary of the change:
#  - The "today" date placeholder fields throughout the deck (slide master,
#    every slide layout and the notes master) are refreshed from "3/10/19"
#    to "4/7/19".
#  - A small group of shapes inside the one real slide (an ER-style
#    "class diagram" box) are tweaked:
#      * the "Name" box is nudged by a single EMU in x/cx (rounding noise
#        from being re-positioned in the UI)
#      * the "Email" box is relabelled "Date"
#      * the "Client" box is relabelled "Serial Number" and widened
#      * a decision/connector/"Address" trio that used to sit to the right
#        of it is removed completely
#      * the trailing "1" callout textbox is removed (the "*" callout
#        right before it is left alone)

$p = $ppt.ActivePresentation

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Id -eq $id) {
            return $shp
        }
    }
    return $null
}

function Update-DateFigureOut($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.TextRange.Text -eq "3/10/19") {
                $shp.TextFrame.TextRange.Text = "4/7/19"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 1. Refresh the "datetimeFigureOut" placeholders (slide master, every
#    layout, and the notes master).
# ---------------------------------------------------------------------

Update-DateFigureOut $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DateFigureOut $layouts.Item($li).Shapes
}

try {
    Update-DateFigureOut $p.NotesMaster.Shapes
} catch {
    # best effort -- not all hosts allow editing the notes master text
}

# ---------------------------------------------------------------------
# 2. Tweak shapes on slide 1.
# ---------------------------------------------------------------------

$s = $p.Slides.Item(1)

# "Name" rectangle -- 1-EMU nudge in x/cx.
$nameBox = Get-ShapeById $s.Shapes 76
$nameBox.Left = 6839239 / 12700.0
$nameBox.Width = 708187 / 12700.0

# "Email" -> "Date"
$emailBox = Get-ShapeById $s.Shapes 83
$emailBox.TextFrame.TextRange.Text = "Date"

# "Client" -> "Serial Number", and widen the box to fit.
$clientBox = Get-ShapeById $s.Shapes 31
$clientBox.Width = 926553 / 12700.0
$clientBox.TextFrame.TextRange.Text = "Serial Number"

# Remove the "Flowchart: Decision 96" (id 32), its "Elbow Connector 32"
# (id 33), and the "Address" rectangle (id 34) that used to sit to the
# right of the "Client"/"Serial Number" box.
(Get-ShapeById $s.Shapes 32).Delete()
(Get-ShapeById $s.Shapes 33).Delete()
(Get-ShapeById $s.Shapes 34).Delete()

# Remove the trailing "1" callout textbox (id 37). The "*" callout
# (id 36) right before it is left untouched.
(Get-ShapeById $s.Shapes 37).Delete()
